# Adapt column header formatting to respective input file names.
# - Rename "<Name>_old" headers (A1:J1) to "<Name>_FV2210"
# - Rename "<Name>_new" headers (L1:U1) to "<Name>_FV2304"
# - Freeze the header row (pane split after row 1)
# - Turn the data range into an actual Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cells ------------------------------------------------
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = [string]$cell.Value2
    if ($current -like "*_old") {
        $cell.Value = ($current -replace "_old$", "_FV2210")
    } elseif ($current -like "*_new") {
        $cell.Value = ($current -replace "_new$", "_FV2304")
    }
}

# --- Freeze header row ----------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- Convert the data range into a formatted Excel Table -----------------
$tableRange = $ws.Range("A1:U72")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)  # xlSrcRange=1, xlYes=1
$listObject.Name = "Table1"

$wb.Save() | Out-Null
